$d = $word.ActiveDocument

# The paragraph currently reads "Version 1." (chars 0-9) followed by the
# _GoBack bookmark (collapsed, sitting right before the paragraph mark).
# Target: "Version 2." where the runs are split as:
#   "Versi" | "on" | <spellEnd/> | " 2" | <bookmarkStart/><bookmarkEnd/> | "."

# 1) Split the "Version" run into "Versi" + "on" by inserting and
#    immediately deleting a temporary bookmark at the split point (position 5).
#    This forces a clean run break without leaving stray run properties.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("__tmpSplit", $splitPoint)
$d.Bookmarks("__tmpSplit").Delete()

# 2) Change "1" to "2" (position 8, the character right after "Version ").
$digit = $d.Range(8, 9)
$digit.Text = "2"

# 3) Remove the trailing "." (now at position 9) from the " 2." run...
$period = $d.Range(9, 10)
$period.Delete()

# 4) ...and re-insert it immediately after the _GoBack bookmark, so the
#    period ends up in its own run placed after the bookmark.
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertAfter(".")
